$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.940.35"
$ws.Range("D3").Value = "1.636.99"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.46"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.620.72"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'0.569"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "'65.27"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "27.946.04"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'231.06"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("D25").Value = "'153.90"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "'6.98"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'15.63"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.111"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").Value = "1.408.56"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.954"
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0170"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.563"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'0.875"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "1.778.36"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'87.98"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "'7.59"
$ws.Range("E51").Value = "  -1.89%  "
